$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.036.47"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "1.822.87"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.36"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4553"
$ws.Range("E7").Value = "  +6.71%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3713"
$ws.Range("E8").Value = "  +1.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07290"
$ws.Range("E9").Value = "  +1.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8593"
$ws.Range("E10").Value = "  -0.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.81"
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("D12").Value = "1.823.50"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.670"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.86"
$ws.Range("E14").Value = "  +4.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.335"
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07108"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008837"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.98"
$ws.Range("E20").Value = "  -0.44%  "
$ws.Range("D21").Value = "27.089.09"
$ws.Range("E21").Value = "  -0.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.180"
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.96"
$ws.Range("E23").Value = "  +0.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.998"
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.62"
$ws.Range("E25").Value = "  -1.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.220"
$ws.Range("E26").Value = "  +5.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.49"
$ws.Range("E27").Value = "  +0.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.270"
$ws.Range("E28").Value = "  +0.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.05"
$ws.Range("E29").Value = "  +0.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08892"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.189"
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7546"
$ws.Range("E32").Value = "  -0.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.946"
$ws.Range("E33").Value = "  +4.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.460"
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("E35").Value = "  -0.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.098"
$ws.Range("E36").Value = "  -1.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01970"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05255"
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5353"
$ws.Range("E39").Value = "  +6.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.176"
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.881"
$ws.Range("E41").Value = "  -0.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1717"
$ws.Range("E42").Value = "  +2.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5244"
$ws.Range("E43").Value = "  +9.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.548"
$ws.Range("E44").Value = "  -0.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.64"
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.964"
$ws.Range("E46").Value = "  +8.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.60"
$ws.Range("E47").Value = "  -0.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.676"
$ws.Range("E48").Value = "  +1.03%  "
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06410"
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.30"
$ws.Range("E51").Value = "  +0.42%  "
